$d = $word.ActiveDocument

$d.Content.Find.Execute("Questions: Introduction to radians", $true, $false, $false, $false, $false, $true, 1, $false, "Questions: Introduction to radians", 2)
$d.Content.Find.Execute("Mark Toner, Ifan Howell-Baines", $true, $false, $false, $false, $false, $true, 1, $false, "Mark Toner, Ifan Howell-Baines", 2)
$d.Content.Find.Execute("Questions relating to the introduction to radians study guide.", $true, $false, $false, $false, $false, $true, 1, $false, "Questions relating to the introduction to radians study guide.", 2)
